$wb = $excel.ActiveWorkbook

# ---- YDS sheet: append new game data to the long space-separated number strings ----
$ydsWs = $wb.Worksheets("YDS")
$ydsWs.Range("B2").Value = "35 3 6 1 2 2 11 4 5 2 0 3 4 4 1 4 5 1 -1 10 -4 1 3 3 1 8 1 3 5 2 4 1 -1 9 2 5 2 7 7 2 4 2 26 9 16 1 -1 1 5 -2 12 3 -1 -5 7 -2 4 3 6 1 2 0 1 17 17 2 2 2 0 -1 7 12 -4 0 0 0 -1 1 4 4 12 0 7 0 1 7 4 -1 4 0 -3 3 -3 4 7 5 15 4 0 8 3 1 24 9 3 2 5 1 21 4 4 5 3 1 1 1 11 -1 6 9 1 3 4 1 21 6 6 3 1 0 0 6 5 2 12 5 9 2 32 6 1 3 10 10 0 13 3 0 -2 15 5 2 4 9 6 3 -2 1 6 -1 -2 0 1 3 -1 3 4 1 8 5 6 5 6 3 7 4 7 5 2 6 10 4 6 5 4 -2 6 3 2 14 -1 6 3 2 0 1 0 1 1 2 2 2 0 16 12 1 5 2 -1 6 2 13 3 3 2 8 4 11 12 2 -2 3 6 5 1 4 11 2 5 1 7 -6 1 2 9 0 5 -1 15 9 13 4 0 5 5 4 18 2 3 1 5 1 -1 16 11 5 1 5 22 17 2 5 10 -1 2 3 0 1 1 1 21 3 2 1 2 14 17 2 1 5 13 -1 7 -1 5 12 2 10 0 4 1 6 2 7 0 3 1 9 0 2 2 2 2 2 1 6 5 7 4 2 -1 4 6 10 19 -1 7 14 2 4 1 2 4 -1 -1 6 64 10 6 0 16 3 3 2 -2 -4 10 17 2 1 5 4 7 3 4 1 22 11 6 -1 4 1 10 5 5 -3 -5 3 -1 2 2"
$ydsWs.Range("C2").Value = "2 15 4 -1 8 3 -1 -1 2 2 1 15 8 4 -3 5 2 5 4 2 8 13 3 6 1 3 7 3 3 -1 17 5 1 4 4 6 10 -3 4 -2 10 17 3 0 0 4 4 14 11 2 1 0 5 0 7 5 3 6 5 3 3 3 3 8 5 2 7 9 5 7 5 -1 1 -2 1 -3 3 11 3 14 6 3 2 5 5 4 21 -1 11 3 4 3 3 6 3 4 6 8 0 1 7 7 3 8 6 1 4 6 4 1 -2 3 2 8 1 1 2 11 1 2 5 3 0 2 1 3 1 2 1 3 7 1 13 3 4 3 2 2 5 0 2 -1 9 4 1 12 3 4 -1 0 0 0 1 3 9 0 2 4 4 21 8 8 10 2 7 1 3 5 5 0 5 1 3 6 3 1 6 -3 19 4 3 6 3 5 3 10 75 3 3 1 9 3 2 28 4 2 2 3 0 2 4 2 0 4 4 4 -2 3 10 2 1 3 9 6 0 3 0 4 10 1 4 4 5 10 24 19 0 3 1 2 5 2 -3 4 5 3 9 1 3 3 5 -2 1 3 -3 3 0 12 2 1 0 2 5 5 0 -1 -1 4 5 2 11 2 5 2 1 4 2 17 3 13 3 68 6 30 5 6 3 7 14 4 2 8 4 3 1 2 1 5 12 9 2 3 1 3 0 4 0 2 2 1 4 4 5 3 3 1 1 5 3 6 3 2 21 3 0 17 5 1"
$ydsWs.Range("B3").Value = "4 7 4 9 8 5 22 6 7 12 8 25 9 7 21 26 9 6 17 16 4 1 6 10 3 8 3 7 11 7 5 5 2 7 -1 2 28 19 3 24 -1 10 10 3 11 8 9 32 3 -1 1 -1 14 7 2 9 6 17 31 13 7 3 10 8 27 7 9 4 4 4 4 11 5 22 9 10 4 5 -5 10 7 4 5 7 17 1 16 15 5 16 7 3 11 4 12 15 9 4 4 21 13 5 10 12 15 1 12 3 21 30 0 9 7 7 15 11 6 9 17 23 4 5 7 16 9 3 3 21 20 1 -3 4 10 13 24 3 10 9 14 22 20 27 0 8 7 3 5 11 3 11 75 14 28 25 3 24 5 4 19 3 15 28 9 11 4 8 4 1 9 20 4 26 3 22 8 46 22 29 28 14 44 2 5 11 33 7 10 4 15 7 14 2 7 17 15 9 2 7 8 12 -1 4 9 41 3 28 3 7 18 2 12 7 11 20 3 12 26 23 4 8 5 1 10 19 17 0 15 3 26 48 11 5 5 17 6 9 9 19 12 6 9 5 18 6 10 10 6 6 15 3 19 2 10 11 4 22 20 16 13 12 38 5 19 5 3 21 10 7 20 9 11 5 41 4 7 7 12 7"
$ydsWs.Range("C3").Value = "17 5 9 18 23 6 36 4 0 30 12 3 -2 15 3 13 27 7 11 9 8 5 3 17 4 5 3 27 5 8 6 6 12 12 9 12 17 11 11 6 13 4 9 7 10 12 7 8 6 13 13 28 3 6 0 27 8 1 4 28 26 11 11 9 7 23 18 11 8 8 10 7 11 7 67 9 16 13 37 7 15 5 -5 27 27 27 3 14 8 11 14 11 18 17 3 23 1 13 6 13 44 20 11 3 11 4 13 0 18 27 20 1 13 12 33 3 2 4 12 24 14 6 9 7 35 12 13 7 -2 6 15 3 7 22 5 17 21 8 7 21 8 8 12 15 6 7 3 4 6 6 16 12 41 6 1 13 5 5 9 11 11 15 14 7 6 18 11 14 24 13 6 2 6 5 33 13 -3 6 21 19 17 6 9 7 12 6 13 11 11 -1 9 6 9 2 11 5 18 3 9 5 18 6 12 1 9 8 9 9 7 0 10 0 6 5 11 4 19 16 17 4 8 10 2 14 8 1 2 3 24 7 14 -3 5 8 14 16 6 11 19 7 11 7 26 5 7 6"

# ---- ST sheet: append new game data + updated week totals ----
$stWs = $wb.Worksheets("ST")
$stWs.Range("B4").Value = "54 65 61 61 62 56 63 64 51 68 26 64 64 67 67 60 63 61 59 64 67 64 67 65 67 66 66 61 68 56 54 54 60 66 65 64 65 65 64 54 67 62 64 61 67 63 56"
$stWs.Range("B5").Value = "0 38 18 21 23 16 0 13 10 24 0 26 22 24 25 22 27 27 22 23 25 19 25 17 25 29 21 22 19 20 12 22 18 28 19 25 28 21 16 17 27 18 20 22 30 16 27"
$stWs.Range("D3").Value = "37 62 32 48 44 35 60 64 52 58 41 68 37 49 44 59 58 58 44 41 42 34 35 46 38 62 55 42 55 66 44 42 60 47 48 36 15 46 45 39 71"
$stWs.Range("D4").Value = "0 18 0 6 0 0 0 25 0 0 0 0 0 10 0 5 12 11 -1 0 -1 0 0 0 0 12 0 0 4 0 0 0 11 0 0 8 0 0 0 0 0"
$stWs.Range("D5").Value = "7 0 -1 14 0 0 8 0 18 12 0 0 11 6 0 13 0 23 17 13 21 22 10 0 27 0 18 15 7 0 0 0 0 0 10 0 0 17 0 15 8 0 0 0 0"
$stWs.Range("B2").Value = 81
$stWs.Range("D2").Value = 41
$stWs.Range("H2").Value = 3
$stWs.Range("I2").Value = 3
$stWs.Range("J2").Value = 214
$stWs.Range("K2").Value = 198
$stWs.Range("L2").Value = 119
$stWs.Range("M2").Value = 88
$stWs.Range("B3").Value = 34

# ---- OFF sheet: updated week totals ----
$offWs = $wb.Worksheets("OFF")
$offWs.Range("C2").Value = 166
$offWs.Range("E2").Value = 9
$offWs.Range("F2").Value = 46
$offWs.Range("G2").Value = 47
$offWs.Range("H2").Value = 7
$offWs.Range("I2").Value = 5
$offWs.Range("J2").Value = 32
$offWs.Range("O2").Value = 7
$offWs.Range("P2").Value = 5
$offWs.Range("C3").Value = 112
$offWs.Range("E3").Value = 28
$offWs.Range("H3").Value = 16
$offWs.Range("L3").Value = 141
$offWs.Range("M3").Value = 99
$offWs.Range("Q3").Value = 352

# ---- DEF sheet: updated week totals ----
$defWs = $wb.Worksheets("DEF")
$defWs.Range("C2").Value = 158
$defWs.Range("E2").Value = 5
$defWs.Range("F2").Value = 46
$defWs.Range("G2").Value = 46
$defWs.Range("I2").Value = 5
$defWs.Range("J2").Value = 24
$defWs.Range("N2").Value = 17
$defWs.Range("O2").Value = 24
$defWs.Range("C3").Value = 121
$defWs.Range("E3").Value = 25
$defWs.Range("F3").Value = 84
$defWs.Range("G3").Value = 29
$defWs.Range("H3").Value = 24
$defWs.Range("I3").Value = 48
$defWs.Range("J3").Value = 39
$defWs.Range("L3").Value = 190
$defWs.Range("M3").Value = 110
$defWs.Range("Q3").Value = 347

# ---- TURNS sheet: updated week totals ----
$turnsWs = $wb.Worksheets("TURNS")
$turnsWs.Range("D3").Value = 7
$turnsWs.Range("E3").Value = 9

# ---- PEN sheet: updated week totals ----
$penWs = $wb.Worksheets("PEN")
$penWs.Range("B2").Value = 11
$penWs.Range("D2").Value = 6
$penWs.Range("B3").Value = 18
